# Insert two new rows before row 265, shifting the existing 265-271 rows down
# to 267-273 (new data is prepended as the two most-recent weekly entries).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A265:A266").EntireRow.Insert()

# Row 265 - new weekly entry
$ws.Cells.Item(265,1).Value = 4
$ws.Cells.Item(265,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(265,3).Value = "Los Lagos"
$ws.Cells.Item(265,4).Value = 44747
$ws.Cells.Item(265,5).Value = 10
$ws.Cells.Item(265,6).Value = 100112017
$ws.Cells.Item(265,7).Value = "Apio"
$ws.Cells.Item(265,8).Value = "Americana (o)"
$ws.Cells.Item(265,9).Value = "Primera"
$ws.Cells.Item(265,10).Value = 25
$ws.Cells.Item(265,11).Value = 11000
$ws.Cells.Item(265,12).Value = 11000
$ws.Cells.Item(265,13).Value = 11000
$ws.Cells.Item(265,14).Value = "`$/docena de matas"
$ws.Cells.Item(265,15).Value = "Región de Coquimbo"
$ws.Cells.Item(265,16).Value = 1833
$ws.Cells.Item(265,17).Value = 6
$ws.Cells.Item(265,18).Value = "Hortaliza"

# Row 266 - new weekly entry
$ws.Cells.Item(266,1).Value = 4
$ws.Cells.Item(266,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(266,3).Value = "Los Lagos"
$ws.Cells.Item(266,4).Value = 44747
$ws.Cells.Item(266,5).Value = 10
$ws.Cells.Item(266,6).Value = 100112017
$ws.Cells.Item(266,7).Value = "Apio"
$ws.Cells.Item(266,8).Value = "Americana (o)"
$ws.Cells.Item(266,9).Value = "Segunda"
$ws.Cells.Item(266,10).Value = 25
$ws.Cells.Item(266,11).Value = 10000
$ws.Cells.Item(266,12).Value = 10000
$ws.Cells.Item(266,13).Value = 10000
$ws.Cells.Item(266,14).Value = "`$/docena de matas"
$ws.Cells.Item(266,15).Value = "Región de Coquimbo"
$ws.Cells.Item(266,16).Value = 1667
$ws.Cells.Item(266,17).Value = 6
$ws.Cells.Item(266,18).Value = "Hortaliza"
